$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ49040422",
    "summ49211876",
    "summ49403896",
    "summ49641176",
    "summ49912485",
    "summ50173247",
    "summ50407456",
    "summ50630545",
    "summ50855219",
    "summ51117779",
    "summ51341804",
    "summ51561859",
    "summ51797192",
    "summ52025819",
    "summ52291884",
    "summ52543423",
    "summ52796610",
    "summ53043030",
    "summ53293859",
    "summ53629142",
    "summ53863738",
    "summ54110346",
    "summ54363474",
    "summ54582939",
    "summ54821000",
    "summ55037635",
    "summ55277683",
    "summ55501771",
    "summ55790783",
    "summ56007279"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

